$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2..17 (A, B-name, C, D, E-bool)
$data = @(
    @(0,  "line1", 7,  9,  $true),
    @(1,  "line2", 9,  8,  $true),
    @(2,  "line3", 8,  10, $true),
    @(3,  "line4", 8,  11, $false),
    @(4,  "line5", 10, 5,  $true),
    @(5,  "line6", 12, 8,  $true),
    @(6,  "line7", 14, 11, $true),
    @(7,  "line8", 16, 9,  $true),
    @(8,  "extr1", 5,  12, $true),
    @(9,  "extr2", 5,  9,  $true),
    @(10, "extr3", 10, 11, $false),
    @(11, "extr4", 7,  8,  $false),
    @(12, "extr5", 9,  11, $false),
    @(13, "extr6", 7,  11, $true),
    @(14, "extr7", 5,  7,  $false),
    @(15, "extr8", 8,  5,  $false)
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row = $row + 1
}

# Rows 16 and 17 are brand new (table grew from 15 to 17 data rows) -
# copy the A-column formatting (bold/border/center-top style) from the
# last pre-existing row so the new rows match the rest of the table.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
